$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A2: numeric 0, styled like the header row (reuse style index 1) ---
$ws.Range("A2").Value = 0
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- B2:F2, H2, K2: plain text values ---
$ws.Range("B2").Value = "Miguel"
$ws.Range("C2").Value = "Angel"
$ws.Range("D2").Value = "Elizondo"
$ws.Range("E2").Value = "Herrera"
$ws.Range("F2").Value = "Posgrado"
$ws.Range("H2").Value = "20/03/2023"
$ws.Range("K2").Value = "C:/Users/MrJua/Desktop/SCI/Fotos/XXMKYX_00.jpeg"

# --- G2: numeric-looking text "140508" must stay text, default style ---
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "140508"
$ws.Range("G2").Style = $ws.Range("B1").Style

# --- I2: date serial styled with a custom date-time number format ---
$ws.Range("I2").Value = 45371
$ws.Range("I2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("I2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- J2: plain number ---
$ws.Range("J2").Value = 5361564

$excel.CutCopyMode = $false
